$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.4458934228143729
$ws.Range("C2").Value = -0.1134040234363881
$ws.Range("D2").Value = 0.6219785646361177

$ws.Range("B3").Value = -0.04237791528613299
$ws.Range("C3").Value = 0.3807011644083477

$ws.Range("B4").Value = 0.7181631419195285

$ws.Range("B5").Value = 1.049089967442578
$ws.Range("C5").Value = 0.8524940945699551
$ws.Range("D5").Value = -0.1034953623938322
$ws.Range("E5").Value = 0.0926409719578045

$ws.Range("B6").Value = 1.499698792741737
$ws.Range("C6").Value = -0.0441541421789331
$ws.Range("D6").Value = -0.3287619841940873

$ws.Range("B7").Value = 0.8297131584907731
$ws.Range("C7").Value = -0.1093333813988623

$ws.Range("B8").Value = 0.1620838256790951

$ws.Range("B9").Value = 0.3211787666311243
$ws.Range("C9").Value = 0.112909106075748
$ws.Range("D9").Value = 0.03134352168291979
$ws.Range("E9").Value = -0.06584948318657535

$ws.Range("B10").Value = 0.2286550815363005
$ws.Range("C10").Value = 0.125655419861192
$ws.Range("D10").Value = -0.1481903203942864

$ws.Range("B11").Value = 0.4490319412376911
$ws.Range("C11").Value = -0.2127353874755186

$ws.Range("B12").Value = 0.1181141524322671

$ws.Range("B13").Value = -0.0001442643046098147
$ws.Range("C13").Value = -0.008000482844254697
$ws.Range("D13").Value = 0.2025442913845202

$ws.Range("B14").Value = 0.1064660537556523
$ws.Range("C14").Value = 0.1601369838950918

$ws.Range("B15").Value = 0.04651757380517277
